# "checkpoint for video style transfer"
#
# This reproduces the content-level changes from the commit:
#   1. The cached text of every "datetimeFigureOut" date field (one on the
#      slide master, one on each of the 11 slide layouts) is bumped from
#      3/19/2020 -> 3/27/2020.
#   2. Three connector shapes on slide 1 are resized / repositioned and two
#      of them get a slightly different bentConnector3 "adj1" adjustment
#      value (the diagram's elbow connectors were nudged/reflowed).
#
# Note: the canonical diff also adds an empty
#   <p:extLst><p:ext uri="{EFAFB233-...}"><p15:sldGuideLst/></p:ext></p:extLst>
# to ppt/presentation.xml. That element just records an (empty) list of
# PowerPoint 2013+ "static" slide guides; this COM host does not expose a
# working Guides object (Presentation.Guides / Master.Guides / Slide.Guides
# all resolve to Nothing here), so there is no object-model call available
# that can reproduce that purely-structural, content-free marker. Everything
# that affects actual slide content/appearance is applied below.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Date placeholder text: slide master + every custom layout.
# ---------------------------------------------------------------------
function Set-DatePlaceholderText($container, $newText) {
    for ($i = 1; $i -le $container.Shapes.Count; $i++) {
        $shp = $container.Shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = $newText
        }
    }
}

$master = $p.SlideMaster
Set-DatePlaceholderText $master "3/27/2020"

for ($L = 1; $L -le $master.CustomLayouts.Count; $L++) {
    $layout = $master.CustomLayouts.Item($L)
    Set-DatePlaceholderText $layout "3/27/2020"
}

# ---------------------------------------------------------------------
# 2) Reflow the three connector shapes on slide 1.
# ---------------------------------------------------------------------
function Get-ShapeById($container, $id) {
    for ($i = 1; $i -le $container.Shapes.Count; $i++) {
        $shp = $container.Shapes.Item($i)
        if ($shp.Id -eq $id) {
            return $shp
        }
    }
    return $null
}

$slide = $p.Slides.Item(1)

# "Connector: Elbow 72" (id 73): ext cx 3600459 -> 2409199 EMU,
# adj1 100000 -> 99973 (per-mille adjustment value).
$cxn1 = Get-ShapeById $slide 73
$cxn1.Width = 189.7007586614173
$cxn1.Adjustments.Item(1) = 0.99973

# "Connector: Elbow 89" (id 90): ext cx 3473088 -> 2242458 EMU,
# adj1 18114 -> 19757.
$cxn2 = Get-ShapeById $slide 90
$cxn2.Width = 176.57154606299213
$cxn2.Adjustments.Item(1) = 0.19757

# "Straight Arrow Connector 92" (id 93): off x 11181715 -> 9951720 EMU.
$cxn3 = Get-ShapeById $slide 93
$cxn3.Left = 783.60005
